$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1040
$ws.Range("I58").Value = 545
$ws.Range("K58").Value = 1635
$ws.Range("M58").Value = -1485

$ws.Range("H74").Value = 5686.7617
$ws.Range("I74").Value = 3832.8333
$ws.Range("K74").Value = 3832.8333
$ws.Range("M74").Value = -2896.8333

$ws.Range("H76").Value = 3741.4167
$ws.Range("I76").Value = 3733.111
$ws.Range("J76").Value = 3766.3333
$ws.Range("K76").Value = 3733.111
$ws.Range("L76").Value = 3766.3333
$ws.Range("M76").Value = -3418.111
$ws.Range("N76").Value = -4396.3333

$ws.Range("H77").Value = 5686.7617
$ws.Range("I77").Value = 3832.8333
$ws.Range("K77").Value = 19164.1665
$ws.Range("M77").Value = -14484.1665

$ws.Range("H79").Value = 3741.4167
$ws.Range("I79").Value = 3733.111
$ws.Range("J79").Value = 3766.3333
$ws.Range("K79").Value = 3733.111
$ws.Range("L79").Value = 3766.3333
$ws.Range("M79").Value = -2641.111
$ws.Range("N79").Value = -5950.3333

$ws.Range("H86").Value = 2638.0527
$ws.Range("I86").Value = 2126
$ws.Range("K86").Value = 2126
$ws.Range("M86").Value = -1003

$ws.Range("H89").Value = 2638.0527
$ws.Range("I89").Value = 2126
$ws.Range("K89").Value = 10630
$ws.Range("M89").Value = -5014

$ws.Range("H138").Value = 2204.318
$ws.Range("I138").Value = 1626.5
$ws.Range("J138").Value = 2782.1365
$ws.Range("K138").Value = 4879.5
$ws.Range("L138").Value = 8346.4095
$ws.Range("M138").Value = 260.5
$ws.Range("N138").Value = -18626.4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6806484
$ws.Range("I32").Value = 7411240
$ws.Range("K32").Value = 7411240
$ws.Range("M32").Value = -7410953

$ws.Range("H45").Value = 1827.6666
$ws.Range("I45").Value = 1899.9333
$ws.Range("K45").Value = 1899.9333
$ws.Range("M45").Value = -1522.9333

$ws.Range("H61").Value = 1151399.8
$ws.Range("I61").Value = 1192306.9
$ws.Range("K61").Value = 1192306.9
$ws.Range("M61").Value = -1192094.9

$ws.Range("H122").Value = 3513.8928
$ws.Range("I122").Value = 3421.2778
$ws.Range("K122").Value = 10263.8334
$ws.Range("M122").Value = -7813.8334

$ws.Range("H136").Value = 1151399.8
$ws.Range("I136").Value = 1192306.9
$ws.Range("K136").Value = 3576920.7
$ws.Range("M136").Value = -3574370.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H86").Value = 1733
$ws.Range("I86").Value = 1949.5
$ws.Range("K86").Value = 1949.5
$ws.Range("M86").Value = -826.5

$ws.Range("H89").Value = 1733
$ws.Range("I89").Value = 1949.5
$ws.Range("K89").Value = 9747.5
$ws.Range("M89").Value = -4131.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 96667.836
$ws.Range("J31").Value = 23659.945
$ws.Range("L31").Value = 23659.945
$ws.Range("N31").Value = -24249.945

$ws.Range("H34").Value = 96667.836
$ws.Range("J34").Value = 23659.945
$ws.Range("L34").Value = 23659.945
$ws.Range("N34").Value = -24063.945

$ws.Range("H117").Value = 40712
$ws.Range("J117").Value = 40712
$ws.Range("L117").Value = 40712
$ws.Range("N117").Value = -49890

$ws.Range("H132").Value = 22847020
$ws.Range("I132").Value = 27029968
$ws.Range("K132").Value = 81089904
$ws.Range("M132").Value = -81087374

$ws.Range("H134").Value = 3416974.2
$ws.Range("I134").Value = 9407.706
$ws.Range("K134").Value = 28223.118
$ws.Range("M134").Value = -25688.118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 222.25
$ws.Range("J5").Value = 199.5
$ws.Range("L5").Value = 598.5
$ws.Range("N5").Value = -822.5

$ws.Range("H57").Value = 1725.5
$ws.Range("I57").Value = 1725.5
$ws.Range("K57").Value = 5176.5
$ws.Range("M57").Value = -4617.5

$ws.Range("H60").Value = 485.08334
$ws.Range("I60").Value = 196.9
$ws.Range("K60").Value = 590.7
$ws.Range("M60").Value = -339.7

$ws.Range("H75").Value = 6618
$ws.Range("I75").Value = 1215.8
$ws.Range("J75").Value = 9994.375
$ws.Range("K75").Value = 3647.4
$ws.Range("L75").Value = 29983.125
$ws.Range("M75").Value = -2649.4
$ws.Range("N75").Value = -31979.125

$ws.Range("H78").Value = 6618
$ws.Range("I78").Value = 1215.8
$ws.Range("J78").Value = 9994.375
$ws.Range("K78").Value = 10942.2
$ws.Range("L78").Value = 89949.375
$ws.Range("M78").Value = -5950.199999999999
$ws.Range("N78").Value = -99933.375

$ws.Range("H98").Value = 1202.3334
$ws.Range("J98").Value = 1402
$ws.Range("L98").Value = 4206
$ws.Range("N98").Value = -7202

$ws.Range("H106").Value = 6338.75
$ws.Range("I106").Value = 7429
$ws.Range("J106").Value = 5248.5
$ws.Range("K106").Value = 22287
$ws.Range("L106").Value = 15745.5
$ws.Range("M106").Value = -21341
$ws.Range("N106").Value = -17637.5

$ws.Range("H107").Value = 790.4
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 800.5
$ws.Range("K107").Value = 2250
$ws.Range("L107").Value = 2401.5
$ws.Range("M107").Value = -330
$ws.Range("N107").Value = -6241.5

$ws.Range("H112").Value = 3465
$ws.Range("J112").Value = 3465
$ws.Range("L112").Value = 10395
$ws.Range("N112").Value = -12611

$ws.Range("H132").Value = 1695.3077
$ws.Range("I132").Value = 919.1
$ws.Range("K132").Value = 8271.9
$ws.Range("M132").Value = -5741.9

$ws.Range("H135").Value = 222.25
$ws.Range("J135").Value = 199.5
$ws.Range("L135").Value = 1795.5
$ws.Range("N135").Value = -6865.5

$ws.Range("H136").Value = 8575
$ws.Range("I136").Value = 7626.6665
$ws.Range("J136").Value = 9997.5
$ws.Range("K136").Value = 22879.9995
$ws.Range("L136").Value = 29992.5
$ws.Range("M136").Value = -17779.9995
$ws.Range("N136").Value = -40192.5

$ws.Range("H138").Value = 6076.143
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 6076.143
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 18228.429
$ws.Range("N138").Value = -28508.429
$ws.Range("M138").ClearContents()

$ws.Range("H139").Value = 74252.78999999999
$ws.Range("I139").Value = 84961.586
$ws.Range("K139").Value = 254884.758
$ws.Range("M139").Value = -249744.758

$ws.Range("H140").Value = 4669.8887
$ws.Range("I140").Value = 7014.5
$ws.Range("K140").Value = 21043.5
$ws.Range("M140").Value = -15863.5

$ws.Range("H141").Value = 3519.6667
$ws.Range("I141").Value = 3519.6667
$ws.Range("K141").Value = 10559.0001
$ws.Range("M141").Value = -5379.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 27415.666
$ws.Range("J96").Value = 27415.666
$ws.Range("L96").Value = 27415.666
$ws.Range("N96").Value = -32907.666

$ws.Range("H122").Value = 40227.38
$ws.Range("J122").Value = 7591.2354
$ws.Range("L122").Value = 22773.7062
$ws.Range("N122").Value = -27673.7062

$ws.Range("H132").Value = 21091108
$ws.Range("I132").Value = 25956274
$ws.Range("K132").Value = 77868822
$ws.Range("M132").Value = -77866292

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 573.4
$ws.Range("I22").Value = 529.5
$ws.Range("J22").Value = 749
$ws.Range("K22").Value = 529.5
$ws.Range("L22").Value = 749
$ws.Range("M22").Value = -234.5
$ws.Range("N22").Value = -1339

$ws.Range("H27").Value = 573.4
$ws.Range("I27").Value = 529.5
$ws.Range("J27").Value = 749
$ws.Range("K27").Value = 529.5
$ws.Range("L27").Value = 749
$ws.Range("M27").Value = -422.5
$ws.Range("N27").Value = -963

$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 3914.1428
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 3914.1428
$ws.Range("M40").Value = -4864
$ws.Range("N40").Value = -4186.1428

$ws.Range("H55").Value = 17857388
$ws.Range("I55").Value = 280
$ws.Range("K55").Value = 280
$ws.Range("M55").Value = -107

$ws.Range("H82").Value = 1014.2
$ws.Range("J82").Value = 1075.4615
$ws.Range("L82").Value = 1075.4615
$ws.Range("N82").Value = -1797.4615

$ws.Range("H85").Value = 1014.2
$ws.Range("J85").Value = 1075.4615
$ws.Range("L85").Value = 1075.4615
$ws.Range("N85").Value = -3571.4615

$ws.Range("H93").Value = 1591.5
$ws.Range("I93").Value = 911.5
$ws.Range("K93").Value = 911.5
$ws.Range("M93").Value = 336.5

$ws.Range("H104").Value = 87958.2
$ws.Range("J104").Value = 87958.2
$ws.Range("L104").Value = 87958.2
$ws.Range("N104").Value = -94946.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 74998
$ws.Range("J123").Value = 74998
$ws.Range("L123").Value = 74998
$ws.Range("N123").Value = -84798

$ws.Range("H126").Value = 935.8
$ws.Range("I126").Value = 935.8
$ws.Range("K126").Value = 2807.4
$ws.Range("M126").Value = -337.3999999999996

$ws.Range("H136").Value = 52366.5
$ws.Range("I136").Value = 59159.5
$ws.Range("K136").Value = 177478.5
$ws.Range("M136").Value = -174928.5
